$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PERSONANATURALGENERAL")
$ws.Activate() | Out-Null

# Update the name data in row 3 (person with PEGE_ID 2222):
# previously placeholder values AAAA/BBBB/CCCC/DDDD, now a real name
$ws.Range("B3").Value = "BOLAÑOS"
$ws.Range("C3").Value = "TAZ"
$ws.Range("D3").Value = "CARLOS"
$ws.Range("E3").Value = "ANDRES"

# Underline the PEGE_ID of row 4 to highlight it
$ws.Range("A4").Font.Underline = $true

# Restore the cursor/selection to B2
$ws.Range("B2").Select() | Out-Null
